# Auto-generated edit script for Famfrit_Profits workbook update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# across ALC, ARM, BSM, CUL, GSM, LTW, WVR sheets per upstream price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 395
$ws.Range("I5").Value = 395
$ws.Range("K5").Value = 395
$ws.Range("M5").Value = -280
$ws.Range("H132").Value = 4421.1924
$ws.Range("I132").Value = 4679.136
$ws.Range("K132").Value = 14037.408
$ws.Range("M132").Value = -11507.408
$ws.Range("H133").Value = 115489
$ws.Range("J133").Value = 115489
$ws.Range("L133").Value = 115489
$ws.Range("N133").Value = -125609
$ws.Range("H138").Value = 4941.1055
$ws.Range("J138").Value = 8358.200000000001
$ws.Range("L138").Value = 25074.6
$ws.Range("N138").Value = -35354.60000000001
$ws.Range("H141").Value = 7018.8
$ws.Range("I141").Value = 7997.25
$ws.Range("K141").Value = 23991.75
$ws.Range("M141").Value = -18811.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2837.4
$ws.Range("I61").Value = 2463.5
$ws.Range("J61").Value = 4333
$ws.Range("K61").Value = 2463.5
$ws.Range("L61").Value = 4333
$ws.Range("M61").Value = -2251.5
$ws.Range("N61").Value = -4757
$ws.Range("H74").Value = 25073.367
$ws.Range("I74").Value = 25073.367
$ws.Range("K74").Value = 25073.367
$ws.Range("M74").Value = -24199.367
$ws.Range("H77").Value = 25073.367
$ws.Range("I77").Value = 25073.367
$ws.Range("K77").Value = 125366.835
$ws.Range("M77").Value = -120998.835
$ws.Range("H132").Value = 41291.51
$ws.Range("I132").Value = 3796.5312
$ws.Range("J132").Value = 281259.4
$ws.Range("K132").Value = 11389.5936
$ws.Range("L132").Value = 843778.2000000001
$ws.Range("M132").Value = -8859.5936
$ws.Range("N132").Value = -848838.2000000001
$ws.Range("H136").Value = 2837.4
$ws.Range("I136").Value = 2463.5
$ws.Range("J136").Value = 4333
$ws.Range("K136").Value = 7390.5
$ws.Range("L136").Value = 12999
$ws.Range("M136").Value = -4840.5
$ws.Range("N136").Value = -18099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3489.5
$ws.Range("I20").Value = 1901
$ws.Range("J20").Value = 5872.25
$ws.Range("K20").Value = 1901
$ws.Range("L20").Value = 5872.25
$ws.Range("M20").Value = -1654
$ws.Range("N20").Value = -6366.25
$ws.Range("H86").Value = 80332.78
$ws.Range("I86").Value = 65499.168
$ws.Range("J86").Value = 110000
$ws.Range("K86").Value = 65499.168
$ws.Range("L86").Value = 110000
$ws.Range("M86").Value = -64376.168
$ws.Range("N86").Value = -112246
$ws.Range("H89").Value = 80332.78
$ws.Range("I89").Value = 65499.168
$ws.Range("J89").Value = 110000
$ws.Range("K89").Value = 327495.84
$ws.Range("L89").Value = 550000
$ws.Range("M89").Value = -321879.84
$ws.Range("N89").Value = -561232
$ws.Range("I99").Value = 2432.4
$ws.Range("J99").Value = 5068.1665
$ws.Range("K99").Value = 2432.4
$ws.Range("L99").Value = 5068.1665
$ws.Range("M99").Value = -934.4000000000001
$ws.Range("N99").Value = -8064.1665
$ws.Range("H134").Value = 1482.2858
$ws.Range("I134").Value = 1482.2858
$ws.Range("K134").Value = 4446.857400000001
$ws.Range("M134").Value = -1911.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 222.11111
$ws.Range("I2").Value = 246.75
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 1480.5
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -1367.5
$ws.Range("N2").Value = -376
$ws.Range("H7").Value = 93.90000000000001
$ws.Range("J7").Value = 124.75
$ws.Range("L7").Value = 374.25
$ws.Range("N7").Value = -598.25
$ws.Range("H23").Value = 987.625
$ws.Range("I23").Value = 420
$ws.Range("J23").Value = 1328.2
$ws.Range("K23").Value = 1260
$ws.Range("L23").Value = 3984.6
$ws.Range("M23").Value = -1025
$ws.Range("N23").Value = -4454.6
$ws.Range("H34").Value = 710
$ws.Range("I34").Value = 117.5
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 352.5
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -268.5
$ws.Range("N34").Value = -4668
$ws.Range("H44").Value = 7331.36
$ws.Range("J44").Value = 3600
$ws.Range("L44").Value = 10800
$ws.Range("N44").Value = -11596
$ws.Range("H102").Value = 4874.25
$ws.Range("J102").Value = 5499.3335
$ws.Range("L102").Value = 16498.0005
$ws.Range("N102").Value = -21366.0005
$ws.Range("H137").Value = 2674.6667
$ws.Range("I137").Value = 2390
$ws.Range("J137").Value = 3030.5
$ws.Range("K137").Value = 7170
$ws.Range("L137").Value = 9091.5
$ws.Range("M137").Value = -2070
$ws.Range("N137").Value = -19291.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2230.2222
$ws.Range("I132").Value = 1655.5294
$ws.Range("K132").Value = 4966.5882
$ws.Range("M132").Value = -2436.5882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2375.8696
$ws.Range("I46").Value = 980.61536
$ws.Range("J46").Value = 4189.7
$ws.Range("K46").Value = 980.61536
$ws.Range("L46").Value = 4189.7
$ws.Range("M46").Value = -792.61536
$ws.Range("N46").Value = -4565.7
$ws.Range("H82").Value = 2860.6667
$ws.Range("I82").Value = 2824.75
$ws.Range("J82").Value = 2889.4
$ws.Range("K82").Value = 2824.75
$ws.Range("L82").Value = 2889.4
$ws.Range("M82").Value = -2463.75
$ws.Range("N82").Value = -3611.4
$ws.Range("H85").Value = 2860.6667
$ws.Range("I85").Value = 2824.75
$ws.Range("J85").Value = 2889.4
$ws.Range("K85").Value = 2824.75
$ws.Range("L85").Value = 2889.4
$ws.Range("M85").Value = -1576.75
$ws.Range("N85").Value = -5385.4
$ws.Range("H122").Value = 3910380.8
$ws.Range("I122").Value = 4078.4583
$ws.Range("K122").Value = 12235.3749
$ws.Range("M122").Value = -9785.374899999999
$ws.Range("H132").Value = 2095.1365
$ws.Range("I132").Value = 1502.75
$ws.Range("J132").Value = 2806
$ws.Range("K132").Value = 4508.25
$ws.Range("L132").Value = 8418
$ws.Range("M132").Value = -1978.25
$ws.Range("N132").Value = -13478
$ws.Range("H133").Value = 59550
$ws.Range("J133").Value = 59550
$ws.Range("L133").Value = 59550
$ws.Range("N133").Value = -64610
$ws.Range("H136").Value = 5911.85
$ws.Range("I136").Value = 2952.375
$ws.Range("J136").Value = 17749.75
$ws.Range("K136").Value = 8857.125
$ws.Range("L136").Value = 53249.25
$ws.Range("M136").Value = -6307.125
$ws.Range("N136").Value = -58349.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H62").Value = 8069.4
$ws.Range("I62").Value = 5948.75
$ws.Range("J62").Value = 9483.166999999999
$ws.Range("K62").Value = 5948.75
$ws.Range("L62").Value = 9483.166999999999
$ws.Range("M62").Value = -5324.75
$ws.Range("N62").Value = -10731.167
$ws.Range("H65").Value = 8069.4
$ws.Range("I65").Value = 5948.75
$ws.Range("J65").Value = 9483.166999999999
$ws.Range("K65").Value = 29743.75
$ws.Range("L65").Value = 47415.835
$ws.Range("M65").Value = -26623.75
$ws.Range("N65").Value = -53655.835
$ws.Range("J86").Value = 59193.75
$ws.Range("L86").Value = 59193.75
$ws.Range("N86").Value = -61439.75
$ws.Range("J89").Value = 59193.75
$ws.Range("L89").Value = 295968.75
$ws.Range("N89").Value = -307200.75
$ws.Range("H126").Value = 3985.923
$ws.Range("I126").Value = 3985.923
$ws.Range("K126").Value = 11957.769
$ws.Range("M126").Value = -9487.769
$ws.Range("H132").Value = 2718.2954
$ws.Range("I132").Value = 2542.2632
$ws.Range("K132").Value = 7626.7896
$ws.Range("M132").Value = -5096.7896
$ws.Range("H133").Value = 80305
$ws.Range("J133").Value = 80305
$ws.Range("L133").Value = 80305
$ws.Range("N133").Value = -90425
